# This workbook has a weekly price log for "Ajo" (garlic) at the
# "Macroferia Regional de Talca" market. A new week's record is added as
# the new first data row of this block (row 299); the previously existing
# rows 299-318 shift down by one (to 300-319), preserving their data and
# formatting as Excel does on a real row insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 299, pushing the existing rows 299:318 down to
# 300:319 (this also grows the sheet dimension to A1:R319 and carries the
# row's number formatting down, matching native Excel behavior).
$ws.Rows.Item(299).Insert()

# Populate the newly inserted row 299 with this week's record.
$ws.Cells.Item(299, 1).Value = 5
$ws.Cells.Item(299, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(299, 3).Value = "Maule"
$ws.Cells.Item(299, 4).Value = 44746
$ws.Cells.Item(299, 5).Value = 7
$ws.Cells.Item(299, 6).Value = 100112003
$ws.Cells.Item(299, 7).Value = "Ajo"
$ws.Cells.Item(299, 8).Value = "Chino"
$ws.Cells.Item(299, 9).Value = "1a (cosecha)"
$ws.Cells.Item(299, 10).Value = 300
$ws.Cells.Item(299, 11).Value = 18000
$ws.Cells.Item(299, 12).Value = 18000
$ws.Cells.Item(299, 13).Value = 18000
$ws.Cells.Item(299, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(299, 15).Value = "China"
$ws.Cells.Item(299, 16).Value = 1800
$ws.Cells.Item(299, 17).Value = 10
$ws.Cells.Item(299, 18).Value = "Hortaliza"
